$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formulas in E3:F5 (shifting the schedule dates earlier)
$ws.Range("F2").Formula = "=DATE(2023,1,15)"
$ws.Range("E3").Formula = "=DATE(2023,1,15)"
$ws.Range("F3").Formula = "=DATE(2023,1,16)"
$ws.Range("E4").Formula = "=DATE(2023,1,16)"
$ws.Range("F4").Formula = "=DATE(2023,1,18)"
$ws.Range("E5").Formula = "=DATE(2023,1,18)"

# Update the current selection to match the target state
$ws.Range("H5").Select()
